# Rafraichissement des frames apres un nème uri_load
# Remove columns E:F (birthDate / birthPlace) and rows 12:18 that only had
# data in those removed columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 12 through 18 (only contained E/F data which is being removed)
$ws.Range("A12:F18").EntireRow.Delete() | Out-Null

# Delete columns E:F entirely
$ws.Range("E1:F11").EntireColumn.Delete() | Out-Null
